# Apply numeric corrections to crafting-leve profit sheets (scheduled runner sync)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1495176.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1495176.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4485529.5
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -4485865.5
$ws.Range("H98").Value = 2000.5625
$ws.Range("I98").Value = 1784.24
$ws.Range("J98").Value = 2773.1428
$ws.Range("K98").Value = 1784.24
$ws.Range("L98").Value = 2773.1428
$ws.Range("M98").Value = -286.24
$ws.Range("N98").Value = -5769.1428
$ws.Range("H118").Value = 1360
$ws.Range("I118").Value = 305.7143
$ws.Range("K118").Value = 917.1428999999999
$ws.Range("M118").Value = 739.8571000000001
$ws.Range("H122").Value = 2000.5625
$ws.Range("I122").Value = 1784.24
$ws.Range("J122").Value = 2773.1428
$ws.Range("K122").Value = 5352.72
$ws.Range("L122").Value = 8319.428400000001
$ws.Range("M122").Value = -2902.72
$ws.Range("N122").Value = -13219.4284
$ws.Range("H138").Value = 2195.15
$ws.Range("I138").Value = 1447.7354
$ws.Range("J138").Value = 3172.5386
$ws.Range("K138").Value = 4343.206200000001
$ws.Range("L138").Value = 9517.6158
$ws.Range("M138").Value = 796.7937999999995
$ws.Range("N138").Value = -19797.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 22257
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 22257
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 22257
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -23723

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1858.9412
$ws.Range("I105").Value = 1327.2727
$ws.Range("J105").Value = 2833.6667
$ws.Range("K105").Value = 1327.2727
$ws.Range("L105").Value = 2833.6667
$ws.Range("M105").Value = 419.7273
$ws.Range("N105").Value = -6327.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3813.1875
$ws.Range("I16").Value = 1901.1
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 1901.1
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = -1614.1
$ws.Range("N16").Value = -7574
$ws.Range("H99").Value = 1684.6666
$ws.Range("I99").Value = 1608
$ws.Range("J99").Value = 1838
$ws.Range("K99").Value = 1608
$ws.Range("L99").Value = 1838
$ws.Range("M99").Value = -110
$ws.Range("N99").Value = -4834
$ws.Range("H113").Value = 3813.1875
$ws.Range("I113").Value = 1901.1
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 1901.1
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = 268.9000000000001
$ws.Range("N113").Value = -11340
$ws.Range("H126").Value = 1684.6666
$ws.Range("I126").Value = 1608
$ws.Range("J126").Value = 1838
$ws.Range("K126").Value = 4824
$ws.Range("L126").Value = 5514
$ws.Range("M126").Value = -2354
$ws.Range("N126").Value = -10454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2456.9524
$ws.Range("I70").Value = 1817.8182
$ws.Range("J70").Value = 3160
$ws.Range("K70").Value = 5453.4546
$ws.Range("L70").Value = 9480
$ws.Range("M70").Value = -5138.4546
$ws.Range("N70").Value = -10110
$ws.Range("H73").Value = 2456.9524
$ws.Range("I73").Value = 1817.8182
$ws.Range("J73").Value = 3160
$ws.Range("K73").Value = 5453.4546
$ws.Range("L73").Value = 9480
$ws.Range("M73").Value = -4361.4546
$ws.Range("N73").Value = -11664
$ws.Range("H75").Value = 1081.4546
$ws.Range("I75").Value = 452.2
$ws.Range("J75").Value = 1605.8334
$ws.Range("K75").Value = 1356.6
$ws.Range("L75").Value = 4817.5002
$ws.Range("M75").Value = -358.5999999999999
$ws.Range("N75").Value = -6813.5002
$ws.Range("H78").Value = 1081.4546
$ws.Range("I78").Value = 452.2
$ws.Range("J78").Value = 1605.8334
$ws.Range("K78").Value = 4069.8
$ws.Range("L78").Value = 14452.5006
$ws.Range("M78").Value = 922.2000000000003
$ws.Range("N78").Value = -24436.5006
$ws.Range("H87").Value = 2702.5
$ws.Range("I87").Value = 2046.5625
$ws.Range("J87").Value = 7950
$ws.Range("K87").Value = 6139.6875
$ws.Range("L87").Value = 23850
$ws.Range("M87").Value = -4891.6875
$ws.Range("N87").Value = -26346
$ws.Range("H90").Value = 2702.5
$ws.Range("I90").Value = 2046.5625
$ws.Range("J90").Value = 7950
$ws.Range("K90").Value = 18419.0625
$ws.Range("L90").Value = 71550
$ws.Range("M90").Value = -12179.0625
$ws.Range("N90").Value = -84030
$ws.Range("H103").Value = 419.125
$ws.Range("I103").Value = 353.25
$ws.Range("K103").Value = 1059.75
$ws.Range("M103").Value = -180.75
$ws.Range("H131").Value = 934.1852
$ws.Range("I131").Value = 345
$ws.Range("J131").Value = 1007.8333
$ws.Range("K131").Value = 1035
$ws.Range("L131").Value = 3023.4999
$ws.Range("M131").Value = 4005
$ws.Range("N131").Value = -13103.4999
$ws.Range("H134").Value = 63596.277
$ws.Range("I134").Value = 75248.87
$ws.Range("J134").Value = 5333.3335
$ws.Range("K134").Value = 225746.61
$ws.Range("L134").Value = 16000.0005
$ws.Range("M134").Value = -220676.61
$ws.Range("N134").Value = -26140.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 19400
$ws.Range("J124").Value = 19400
$ws.Range("L124").Value = 19400
$ws.Range("N124").Value = -29220

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1512.3572
$ws.Range("I82").Value = 1297.2
$ws.Range("J82").Value = 1631.8889
$ws.Range("K82").Value = 1297.2
$ws.Range("L82").Value = 1631.8889
$ws.Range("M82").Value = -936.2
$ws.Range("N82").Value = -2353.8889
$ws.Range("H85").Value = 1512.3572
$ws.Range("I85").Value = 1297.2
$ws.Range("J85").Value = 1631.8889
$ws.Range("K85").Value = 1297.2
$ws.Range("L85").Value = 1631.8889
$ws.Range("M85").Value = -49.20000000000005
$ws.Range("N85").Value = -4127.8889
$ws.Range("H132").Value = 18512.297
$ws.Range("I132").Value = 9008.764999999999
$ws.Range("J132").Value = 34668.3
$ws.Range("K132").Value = 27026.295
$ws.Range("L132").Value = 104004.9
$ws.Range("M132").Value = -24496.295
$ws.Range("N132").Value = -109064.9
$ws.Range("H134").Value = 39429
$ws.Range("J134").Value = 39429
$ws.Range("L134").Value = 39429
$ws.Range("N134").Value = -49569
$ws.Range("H136").Value = 4919.2354
$ws.Range("I136").Value = 1338.7727
$ws.Range("J136").Value = 11483.417
$ws.Range("K136").Value = 4016.3181
$ws.Range("L136").Value = 34450.251
$ws.Range("M136").Value = -1466.3181
$ws.Range("N136").Value = -39550.251

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 968.1053000000001
$ws.Range("I122").Value = 875.6923
$ws.Range("J122").Value = 1168.3334
$ws.Range("K122").Value = 2627.0769
$ws.Range("L122").Value = 3505.0002
$ws.Range("M122").Value = -177.0769
$ws.Range("N122").Value = -8405.0002
$ws.Range("H126").Value = 1583.3334
$ws.Range("I126").Value = 1647
$ws.Range("K126").Value = 4941
$ws.Range("M126").Value = -2471
$ws.Range("H136").Value = 5129.625
$ws.Range("I136").Value = 8157.769
$ws.Range("J136").Value = 1550.909
$ws.Range("K136").Value = 24473.307
$ws.Range("L136").Value = 4652.727000000001
$ws.Range("M136").Value = -21923.307
$ws.Range("N136").Value = -9752.727000000001

